$wb = $excel.ActiveWorkbook

# Update the Date value on the Metadata sheet
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-05-05T11:54:16+00:00"

# Update Min/Max and Base Min/Base Max for the ActorPatient.XCN9.composant1 row
# on the Elements sheet (row 7: F=Min, G=Max, AG=Base Min, AH=Base Max) from "1" to "0".
# F2 already holds the text value "0" with the same style as F7/G7/AG7/AH7, so copy
# it across instead of re-typing the literal (which Excel would auto-coerce to a
# number and stamp with a "stored as text" quote-prefix style).
$elements = $wb.Worksheets.Item("Elements")

$elements.Range("F2").Copy()
$elements.Range("F7").PasteSpecial()

$elements.Range("F2").Copy()
$elements.Range("G7").PasteSpecial()

$elements.Range("F2").Copy()
$elements.Range("AG7").PasteSpecial()

$elements.Range("F2").Copy()
$elements.Range("AH7").PasteSpecial()
